$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 265, shifting existing rows 265:320 down to 266:321
$ws.Rows.Item(265).Insert()

# Populate the new row 265 with the new data point
$ws.Cells.Item(265, 1).Value = 3
$ws.Cells.Item(265, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(265, 3).Value = "Coquimbo"
$ws.Cells.Item(265, 4).Value = 44711
$ws.Cells.Item(265, 5).Value = 5
$ws.Cells.Item(265, 6).Value = 100112039
$ws.Cells.Item(265, 7).Value = "Ciboulette"
$ws.Cells.Item(265, 8).Value = "Sin especificar"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 260
$ws.Cells.Item(265, 11).Value = 1500
$ws.Cells.Item(265, 12).Value = 1500
$ws.Cells.Item(265, 13).Value = 1500
$ws.Cells.Item(265, 14).Value = "`$/docena de atados"
$ws.Cells.Item(265, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(265, 16).Value = 500
$ws.Cells.Item(265, 17).Value = 3
$ws.Cells.Item(265, 18).Value = "Hortaliza"
